$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update header row
$ws.Range("B1").Value = "Host"
$ws.Range("G1").Value = "Allow FS Deletion (Yes / No)"
$ws.Range("T1").Value = "Exclude File System(s)"
$ws.Range("U1").Value = "Include File System(s)"

# Row 2: fill in NA for V2 and W2
$ws.Range("V2").Value = "NA"
$ws.Range("W2").Value = "NA"
$ws.Range("V2:W2").Style = "Normal"
$ws.Range("V2:W2").Font.Bold = $false

# Row 3: change B3 text, add V3/W3
$ws.Range("B3").Value = "psp-MyLinSecondFlow-src1, psp-MyLinSecondFlow-src2"
$ws.Range("V3").Value = "NA"
$ws.Range("W3").Value = "NA"

# Row 4: new row
$ws.Range("A4").Value = "Second Flow"
$ws.Range("B4").Value = "psp-MyWinSecondFlow-src1, psp-MyWinSecondFlow-src2"
$ws.Range("T4").Value = "E:"
$ws.Range("U4").Value = "NA"
$ws.Range("V4").Value = "NA"
$ws.Range("W4").Value = "NA"

# Row 5: new row
$ws.Range("A5").Value = "Second Flow"
$ws.Range("H5").Value = "No"
$ws.Range("T5").Value = "NA"
$ws.Range("U5").Value = "NA"
$ws.Range("V5").Value = "NA"
$ws.Range("W5").Value = "NA"

# Column widths (approximate auto-fit result for new content)
$ws.Range("A1").ColumnWidth = 12
$ws.Range("B1").ColumnWidth = 48

# Sheet view adjustments: move selection to A5 (also clears the stale topLeftCell/old selection)
$ws.Range("A5").Select()
